$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title,
#    before the "Gameplay Mechanics" heading.
# ---------------------------------------------------------------------------
$titlePar = $d.Paragraphs.Item(1)
$titlePar.Range.InsertParagraphAfter()

$metaPar = $d.Paragraphs.Item(2)
$metaPar.Style = "Normal"

$metaBoldLabel = "Meta description"
$metaRest = ": Discover the unique game mechanics and bonus features of Blazin Hot 7s Stack Em Up. Play for free and win big with high volatility and RTP."

# Write the full text first, then apply Bold only to the "Meta description"
# part, leaving the remainder (starting with the colon) at normal weight.
$metaPar.Range.Text = $metaBoldLabel + $metaRest

$metaPar = $d.Paragraphs.Item(2)
$boldRange = $d.Range($metaPar.Range.Start, $metaPar.Range.Start + $metaBoldLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2) Remove the duplicate title paragraph near the end of the document and
#    rewrite the following italic paragraph's text with the new image prompt.
# ---------------------------------------------------------------------------
$dupTitleText = "Play Blazin Hot 7s Stack Em Up for Free - Game Review"

$count = $d.Paragraphs.Count
$dupPar = $null
for ($i = $count; $i -ge 1; $i--) {
    $para = $d.Paragraphs.Item($i)
    $paraText = $para.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq $dupTitleText) {
        $dupPar = $para
        break
    }
}
$dupPar.Range.Delete()

# The italic "description" paragraph is now the last paragraph in the
# document; replace its text while keeping italic formatting, using
# InsertAfter (rather than setting .Text) so smart-quote autocorrect does
# not mangle the straight quotes/apostrophe in the new copy.
$newDescText = 'Create a feature image fitting "Blazin Hot 7s Stack Em Up": - Draw a cartoon-style image of a happy Maya warrior with glasses wearing a headdress made of fruits such as cherries, oranges, lemons, plums, and watermelons. - Have the warrior holding a Stack''Em Up symbol in one hand and a handful of coins in the other hand. - Surround the warrior with cascading reels and colorful symbols. - Add text above the image that says "Blazin Hot 7s Stack Em Up" in bold, fiery letters.'

$descPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $descPar.Range
$textRange = $d.Range($r.Start, $r.End - 1)
$textRange.Delete()

$descPar = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPos = $descPar.Range.Start
$ins = $d.Range($insertPos, $insertPos)
$ins.InsertAfter($newDescText)
$ins.Italic = 1

Write-Output "done"
